$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# "Generate Report for Handoff": the 0822cacf... source file is being
# handed off again (fresh handoff timestamp, status flips from
# "Handed back: in sync with en-US" back to "Ready for handoff"), and the
# now-stale a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md row is removed from the
# report entirely (its localization cycle is over).
# -----------------------------------------------------------------------

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Row 3 on every sheet is the a19ccab5-e9ca-4cab-ad56-f8ee072f1184.md entry;
# deleting it shifts the .localization-config row up to row 3.
$overview.Rows.Item(3).Delete()
$zhcn.Rows.Item(3).Delete()
$dede.Rows.Item(3).Delete()

# Excel does not renumber the leftover Hyperlinks collection when rows
# shift, so rebuild it from scratch for the two rows that remain.
$overview.Range("A1:C3").Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c982aa621d06db827e4f56f06a10914408aad719/e2e/0822cacf-b845-43cd-b6ac-8d79fdd175df.md", "", "", "0822cacf-b845-43cd-b6ac-8d79fdd175df.md")
$overview.Hyperlinks.Add($overview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c982aa621d06db827e4f56f06a10914408aad719/.localization-config", "", "", ".localization-config")
$overview.Range("A2").Style = "HyperLink"
$overview.Range("A3").Style = "HyperLink"

$zhcn.Range("A1:I3").Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c982aa621d06db827e4f56f06a10914408aad719/e2e/0822cacf-b845-43cd-b6ac-8d79fdd175df.md", "", "", "0822cacf-b845-43cd-b6ac-8d79fdd175df.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7cec59ea5ea2c2f12ca37fd5ded16a63aab5762/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.zh-cn.xlf", "", "", "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d59395d057b8dbca0ee35c26e8c736469d787aaf/e2e/0822cacf-b845-43cd-b6ac-8d79fdd175df.md", "", "", "0822cacf-b845-43cd-b6ac-8d79fdd175df.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/08bb596f072d0fb1836184f59f099c5bb2edf772/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.zh-cn.xlf", "", "", "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c982aa621d06db827e4f56f06a10914408aad719/.localization-config", "", "", ".localization-config")
$zhcn.Range("A2").Style = "HyperLink"
$zhcn.Range("C2").Style = "HyperLink"
$zhcn.Range("E2").Style = "HyperLink"
$zhcn.Range("F2").Style = "HyperLink"
$zhcn.Range("A3").Style = "HyperLink"

$dede.Range("A1:I3").Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c982aa621d06db827e4f56f06a10914408aad719/e2e/0822cacf-b845-43cd-b6ac-8d79fdd175df.md", "", "", "0822cacf-b845-43cd-b6ac-8d79fdd175df.md")
$dede.Hyperlinks.Add($dede.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a73b4b6a4f03e481a542db6c770188779c9d20b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.de-de.xlf", "", "", "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/429bfdbea27b2637ccd450a77c5d5abd7a74716e/e2e/0822cacf-b845-43cd-b6ac-8d79fdd175df.md", "", "", "0822cacf-b845-43cd-b6ac-8d79fdd175df.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e516b867a9f2d0a6daf4b69b9bf33334ab0dac8c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.de-de.xlf", "", "", "0822cacf-b845-43cd-b6ac-8d79fdd175df.c77dcca1746c842944ac1bea1dd5679f0a71b7d6.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c982aa621d06db827e4f56f06a10914408aad719/.localization-config", "", "", ".localization-config")
$dede.Range("A2").Style = "HyperLink"
$dede.Range("C2").Style = "HyperLink"
$dede.Range("E2").Style = "HyperLink"
$dede.Range("F2").Style = "HyperLink"
$dede.Range("A3").Style = "HyperLink"

# Flip the 0822cacf... row back to "Ready for handoff" with a fresh
# handoff timestamp on each localized sheet.
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

$zhcn.Range("B2").Value = "Ready for handoff"
$zhcn.Range("D2").Value = "2016-03-08 12:35:09"

$dede.Range("B2").Value = "Ready for handoff"
$dede.Range("D2").Value = "2016-03-08 12:35:14"
